# Daily attendance processing - 2025-12-06 05:50:10
#
# Column G ("Recorded By") holds a comma-separated list of the users who
# recorded/touched a given attendance session. Re-processing the log moves
# any "System" marker to the end of its list (System is always the backup
# recorder, so it should be reported last), while keeping the other
# recorder names in their existing relative order. When a list has no
# "System" entry at all, the two recorder names are reported in reverse
# order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $value = $cell.Value2

    if ($value -eq $null -or $value -eq "") {
        continue
    }

    $parts = @($value -split "," | ForEach-Object { $_.Trim() })

    $systemParts = @($parts | Where-Object { $_.ToLower() -eq "system" })
    $otherParts  = @($parts | Where-Object { $_.ToLower() -ne "system" })

    if ($systemParts.Count -gt 0) {
        $newParts = $otherParts + $systemParts
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $cell.Value = [string]::Join(", ", $newParts)
}
